# feat: added seed funders and airdrop allocations
#
# Applies the "raised" sheet restructuring: inserts a new "Token price USD"
# column, introduces SNS Sale / Treasury percentage assumptions (rows 6-7),
# updates the Min/Max raise figures, and repoints the dependent charts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raised")

# --- 1. Make room for the new "Token price USD" column (old E shifts to F,
#        old F shifts to G, etc.) -------------------------------------------
$ws.Columns("E:E").Insert()

# --- 2. Raw input values ----------------------------------------------------
$ws.Range("C2").Value = 100000
$ws.Range("C3").Value = 100000
$ws.Range("B5").Value = 7

$ws.Range("A6").Value = "SNS Sale"
$ws.Range("B6").Value = 0.28000000000000003
$ws.Range("A7").Value = "Treasury"
$ws.Range("B7").Value = 0.43

# --- 3. New column header ----------------------------------------------------
$ws.Range("E1").Value = "Token price USD"

# --- 4. Formulas -------------------------------------------------------------
$ws.Range("D2").Formula = "=(B2+C2)/B6/H2"
$ws.Range("D3").Formula = "=(B3+C3)/B6/H3"

$ws.Range("E2").Formula = "=D2*B5"
$ws.Range("E3").Formula = "=D3*B5"

$ws.Range("F2").Formula = "=J2"
$ws.Range("F3").Formula = "=J3"

$ws.Range("G2").Formula = "=F2*(1-B7)"
$ws.Range("G3").Formula = "=F3*(1-B7)"

$ws.Range("I2").Formula = "=H2*(1-B7)"
$ws.Range("I3").Formula = "=H3*(1-B7)"

# K2:M3, L2:L3 already carry the correct (shifted) formulas after the column
# insert above - only their computed values change, which recalculates
# automatically.

# --- 5. Number formats for the new cells -------------------------------------
# E2:E3 -> same USD currency format used by the M-cap columns (style s=4)
$ws.Range("L2").Copy()
$ws.Range("E2:E3").PasteSpecial(-4122)

# B6:B7 -> same percent format used on the distribution sheet (style s=1)
$wb.Worksheets.Item("distribution").Range("C2").Copy()
$ws.Range("B6:B7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 6. Column widths ---------------------------------------------------------
$ws.Columns("D:E").ColumnWidth = 16.1666666667
$ws.Columns("F:G").ColumnWidth = 15.5924479167

# --- 7. Chart series now point one column further to the right --------------
$charts = $ws.ChartObjects()
for ($i = 1; $i -le $charts.Count; $i++) {
  $chart = $charts.Item($i).Chart
  $series = $chart.SeriesCollection()
  for ($j = 1; $j -le $series.Count; $j++) {
    $ser = $series.Item($j)
    $f = $ser.Formula
    if ($f -like "*raised!`$E*") {
      $ser.Formula = $f.Replace("raised!`$E", "raised!`$F")
    } elseif ($f -like "*raised!`$F*") {
      $ser.Formula = $f.Replace("raised!`$F", "raised!`$G")
    }
  }
}

# --- 8. Selection / view state -----------------------------------------------
$ws.Activate()
$ws.Range("J3").Select()

Write-Host "edit applied"
